$wb = $excel.ActiveWorkbook

# Sheets A1..A5: the second row (generic example/instruction row) is removed,
# shifting the data row up. The user selected the whole row before deleting it,
# which leaves the selection on the new row 2 (A2:XFD2) after the delete.
$simpleSheets = @("A1", "A2", "A3", "A4", "A5")
foreach ($name in $simpleSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows(2).Select()
    $ws.Rows(2).Delete()
}

# Sheet A6 gets the same row-2 deletion, but ends up as the final active sheet
# with a specific cell (C20) selected/clicked afterwards.
$ws6 = $wb.Worksheets.Item("A6")
$ws6.Rows(2).Select()
$ws6.Rows(2).Delete()
$ws6.Range("C20").Select()
$ws6.Activate()
